$d = $word.ActiveDocument

# Locate the paragraph that holds the standalone "…" character -- the new
# "Version management" paragraph goes directly after it, and the new
# "~Cherry1315" signature paragraph goes directly after the following
# (already-existing) empty paragraph.
$ellipsisIndex = 0
$emptyIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "…`r") {
        $ellipsisIndex = $i
    }
    if ($ellipsisIndex -gt 0 -and $t -eq "`r" -and $emptyIndex -eq 0 -and $i -gt $ellipsisIndex) {
        $emptyIndex = $i
    }
}

# --- Insert the "Version management…" paragraph right after the "…" one ---
$ellipsisPara = $d.Paragraphs.Item($ellipsisIndex)
$ellipsisPara.Range.InsertParagraphAfter()
$versionPara = $d.Paragraphs.Item($ellipsisIndex + 1)
$versionPara.Range.Text = "Version management is needed to keep track of the system, especially when multiple people are working on it, and especially if the system needs to revert to an older version due to issues."

# The empty paragraph shifted down by one.
$emptyIndex = $emptyIndex + 1

# --- Insert the "~Cherry1315" signature paragraph right after the empty one ---
$emptyPara = $d.Paragraphs.Item($emptyIndex)
$emptyPara.Range.InsertParagraphAfter()
$signaturePara = $d.Paragraphs.Item($emptyIndex + 1)
$signaturePara.Range.Text = "~"

# Append "Cherry1315" as its own run (distinct from the "~" run) by briefly
# toggling a character property on just that span -- this forces Word to
# keep it as a separate <w:r> rather than coalescing it with the "~" run.
$fullRange = $d.Paragraphs.Item($emptyIndex + 1).Range
$fullRange.Collapse(0)
$fullRange.InsertAfter("Cherry1315")

$nameRange = $d.Paragraphs.Item($emptyIndex + 1).Range
$nameRange.SetRange($nameRange.Start + 1, $nameRange.End - 1)
$nameRange.Bold = 1
$nameRange.Bold = 0
